$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add header for new "Save" column, matching the style of the other header cells (G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the "Save" column values for rows 2-12
$saveValues = @(1, 0, 0, 0, 0, 0, 1, 1, 1, 1, 0)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
